$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column L (2019) into the new column M (2020) for
# every data row, then overwrite the values with the 2020 figures.
$ws.Range("L2:L15").Copy($ws.Range("M2:M15"))

$ws.Range("M3").Value = 2020
$ws.Range("M4").Value = 94.1
$ws.Range("M5").Value = 99.6
$ws.Range("M6").Value = 91
$ws.Range("M7").Value = 86.886172668979881
$ws.Range("M8").Value = 86.955790296225956
$ws.Range("M9").Value = 96.29195112324031
$ws.Range("M10").Value = 97.849780305474511
$ws.Range("M11").Value = 90.676703333930902
$ws.Range("M12").Value = 99.675929342188979
$ws.Range("M13").Value = 100
$ws.Range("M14").Value = 100
$ws.Range("M15").Value = 100

# M4 carries a distinct (new) bold style, unlike the rest of column M which
# reuses the same formatting as the matching cell in column L.
$ws.Range("M4").Font.Bold = $true

# Update the view: select G15, matching the saved workbook state.
$ws.Range("G15").Select()
